$wb = $excel.ActiveWorkbook

# --- Sheet "Sprint 2 Inhalt" (sheet1): scroll down so row 10 area is visible ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
[void]$ws1.Range("D10").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1

# --- Sheet "Sprint 2 Backlog" (sheet2): add new "Konfiguration" topic with a
#     "Jenkins einrichten" subtask, pushing the former row 8 down to row 10 ---
$ws2 = $wb.Worksheets.Item(2)

# Remember the old row-8 values before we overwrite anything.
$oldB8 = $ws2.Range("B8").Value()
$oldC8 = $ws2.Range("C8").Value()
$oldD8 = $ws2.Range("D8").Value()
$oldE8 = $ws2.Range("E8").Value()

# New group header in row 9.
$ws2.Range("A9").Value = "Konfiguration"

# Move the old row 8 subtask down to row 10.
$ws2.Range("B10").Value = $oldB8
$ws2.Range("C10").Value = $oldC8
$ws2.Range("D10").Value = $oldD8
$ws2.Range("E10").Value = $oldE8

# Clear out the now-empty old row 8.
[void]$ws2.Range("B8:E8").ClearContents()

# This sheet becomes the active tab with a new selection.
$ws2.Activate()
[void]$ws2.Range("A16").Select()

# --- Sheet "Product-Backlog" (sheet3): left untouched (no selection/scroll change) ---
